$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the teaching position text for the Stirling row (A10): the old
# "Adistente de Enseñanza" is replaced with the correct title.
$ws.Range("A10").Value = "Profesor Homologado de Psicología"

# Update the selected cell in the sheet view from C20 to A10.
$ws.Range("A10").Select()
